# testing reporting.py for logs feature.
# Append four new log rows (17-20) to the "Sheet" worksheet and drop the
# now-stale SUM formula in "Monthly_STAT"!B2.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet")

# Row 17 - UID + count only
$ws.Cells.Item(17, 1).Value = "saergetsrt"
$ws.Cells.Item(17, 4).Value = 11

# Row 18 - UID, Name, Email, count
$ws.Cells.Item(18, 1).Value = "saeccffrgetsrt"
$ws.Cells.Item(18, 2).Value = "gary"
$ws.Cells.Item(18, 3).Value = "@live"
$ws.Cells.Item(18, 4).Value = 2

# Row 19 - UID, Name, Email, count
$ws.Cells.Item(19, 1).Value = "saeccffrgeaa"
$ws.Cells.Item(19, 2).Value = "Gary Tsai"
$ws.Cells.Item(19, 3).Value = "yue.tsai@jjay.cuny.edu"
$ws.Cells.Item(19, 4).Value = 1

# Row 20 - UID + count only
$ws.Cells.Item(20, 1).Value = "saeccffrgea"
$ws.Cells.Item(20, 4).Value = 3

# The Monthly_STAT summary formula referenced a fixed range that no longer
# covers the new rows correctly, so it was removed.
$ws2 = $wb.Worksheets.Item("Monthly_STAT")
$ws2.Cells.Item(2, 2).ClearContents()

# Leave the selection where the new data entry ended.
[void]$ws.Range("C17").Select()
